# [feature/ResolveReport] email service added
#
# Two previously "reported" appreciation entries get moderated by
# Sharyu Marwadi and their status flips to "deleted", with a moderator
# comment recorded (simulating the email notification service that now
# fires on resolve/delete of a reported appreciation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReportedAppreciations")
$ws.Activate()

# Row 3: Great teamwork aditya! appreciation - reported -> deleted
$ws.Range("N3").Value = "testing email"
$ws.Range("O3").Value = "Sharyu"
$ws.Range("P3").Value = "Marwadi"
$ws.Range("Q3").Value = "deleted"

# Row 13: Great teamwork! appreciation (Arjun/B moderator) - reported -> deleted
$ws.Range("N13").Value = "testing email 2"
$ws.Range("O13").Value = "Sharyu"
$ws.Range("P13").Value = "Marwadi"
$ws.Range("Q13").Value = "deleted"
